# Updated cryptos list on Tue May  7 19:25:56 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the dogwifhat / Stacks rows (37 <-> 38) to reflect the new ranking.
#
# Price values are plain text in this sheet (the original file stores them
# as inlineStr, using "." as a thousands separator in some rows, e.g.
# "63.158.47"), so numeric-looking prices are forced back to text
# (NumberFormat "@") before/after the write and the style is reset to
# "Normal" afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.129.50"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.050.92"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "3.050.20"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "3.550.22"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "63.067.97"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "3.046.95"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.11%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "0.0₃0818"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0362"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "2.826.82"
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -1.66%  "
